$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 21:09"

# 2) In-place statistic refreshes (country keeps its row, only B..H numbers change)
# Row 4: Estados Unidos
$ws.Range("B4").Value = 6310970
$ws.Range("C4").Value = 20233
$ws.Range("D4").Value = 3555933
$ws.Range("E4").Value = 2564580
$ws.Range("G4").Value = 493
$ws.Range("H4").Value = 190457

# Row 6: India
$ws.Range("B6").Value = 3933124
$ws.Range("C6").Value = 84156
$ws.Range("D6").Value = 3032916
$ws.Range("E6").Value = 831639
$ws.Range("G6").Value = 1083
$ws.Range("H6").Value = 68569

# Row 23: Alemania
$ws.Range("B23").Value = 248722
$ws.Range("C23").Value = 1331
$ws.Range("E23").Value = 16224
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 9398

# Row 27: Canada
$ws.Range("B27").Value = 130262
$ws.Range("C27").Value = 339
$ws.Range("D27").Value = 115284
$ws.Range("E27").Value = 5838
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 9140

# Row 80: Libano
$ws.Range("D80").Value = 5338
$ws.Range("E80").Value = 13446

# Row 97: Guayana Francesa
$ws.Range("D97").Value = 8772
$ws.Range("E97").Value = 418

# 3) Re-sorted block: Cuba overtakes Mozambique/Tunez/Eslovaquia/Surinam (rows 118-122)
$ws.Range("A118").Value = 'Cuba'
$ws.Range("B118").Value = 4214
$ws.Range("C118").Value = 88
$ws.Range("D118").Value = 3474
$ws.Range("E118").Value = 640
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 100

$ws.Range("A119").Value = 'Mozambique'
$ws.Range("B119").Value = 4207
$ws.Range("C119").Value = 90
$ws.Range("D119").Value = 2370
$ws.Range("E119").Value = 1811
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 26

$ws.Range("A120").Value = 'Tunez'
$ws.Range("B120").Value = 4196
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 1628
$ws.Range("E120").Value = 2487
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 81

$ws.Range("A121").Value = 'Eslovaquia'
$ws.Range("B121").Value = 4163
$ws.Range("C121").Value = 121
$ws.Range("D121").Value = 2617
$ws.Range("E121").Value = 1509
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 4
$ws.Range("H121").Value = 37

$ws.Range("A122").Value = 'Surinam'
$ws.Range("B122").Value = 4149
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 3272
$ws.Range("E122").Value = 805
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 72

# 4) Siria own update, then Angola overtakes Mali (rows 132, 134-135)
# Row 132: Siria (values refreshed in place)
$ws.Range("B132").Value = 2973
$ws.Range("C132").Value = 75
$ws.Range("D132").Value = 681
$ws.Range("E132").Value = 2168
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 4
$ws.Range("H132").Value = 124

$ws.Range("A134").Value = 'Angola'
$ws.Range("B134").Value = 2805
$ws.Range("C134").Value = 28
$ws.Range("D134").Value = 1144
$ws.Range("E134").Value = 1548
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 113

$ws.Range("A135").Value = 'Mali'
$ws.Range("B135").Value = 2802
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 2185
$ws.Range("E135").Value = 491
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 126

# 5) Re-sorted block: Bonaire overtakes San Bartolome/Islas Virgenes EEUU/San Cristobal (rows 209-212)
$ws.Range("A209").Value = 'Bonaire, San Eustaquio y Saba'
$ws.Range("B209").Value = 18
$ws.Range("C209").Value = 2
$ws.Range("D209").Value = 7
$ws.Range("E209").Value = 11
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = 'San Bartolome'
$ws.Range("B210").Value = 18
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 5
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = 'Islas Virgenes de los Estados Unidos'
$ws.Range("B211").Value = 17
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 0
$ws.Range("E211").Value = 17
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

$ws.Range("A212").Value = 'San Cristobal y Nieves'
$ws.Range("B212").Value = 17
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 17
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

